# Actualización automática 2025-10-13 17:30:09
#
# Updates linked summary values across three worksheets to reflect an
# updated sale amount of 5872.12 (previously 195.64) for
# TOSCANO RAMIREZ MONICA CECILIA under RIOS CARRION ANGEL BENIGNO /
# PORCELANATO, and propagates the resulting totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": detail by client/group ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M22").Value = 5872.12

# --- Sheet "VENTA MENSUAL": detail by client + totals row ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F22").Value = 5872.12
$wsVentaMensual.Range("F26").Value = 15349.44

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row + TOTAL row ---
$wsCumplimientoMensual = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimientoMensual.Range("D12").Value = 16741.47
$wsCumplimientoMensual.Range("E12").Value = 11213.51
$wsCumplimientoMensual.Range("F12").Value = 0.5988725443552455

$wsCumplimientoMensual.Range("D14").Value = 15349.44
$wsCumplimientoMensual.Range("E14").Value = 26853.94110009469
$wsCumplimientoMensual.Range("F14").Value = 0.3637016655986731
